$wb = $excel.ActiveWorkbook

# New "Kadastro" record (Kayit No 2945) appended as the new last row of
# both the master log sheet ("Kayitlar") and the per-birim sheet
# ("Erdemli"). Columns: A=Kayit No, B=Tarih, C=Birim, D=Parsel Sayisi,
# E=Is, F=Personeller.
$kayitNo  = "2945"
$tarih    = "2025-09-08"
$birim    = "Erdemli"
$parsel   = "1"
$is       = "ÇAP"
$personel = "CEMAL TİMUROĞLU (K.Teknisyeni)"

# Helper: writes $text into $cell as a genuine text value (never an
# auto-inferred number/date), without leaving the cell with a different
# style/number-format than its neighbours.
#
# A plain `$cell.Value = "2945"` gets read back by Excel as the *number*
# 2945 (and "2025-09-08" as a date serial) whenever the cell is still
# General-formatted, which is exactly the case here. Forcing the cell to
# Text via NumberFormat = "@" avoids that, but stamps a brand-new style
# onto the cell. Instead, stage the literal text through a formula
# ( ="2945" ) in a scratch cell - its cached result is plain text - then
# copy/paste just that value into the destination, which keeps the
# destination's original (unstyled) look while keeping the text type.
function Set-TextValue {
    param($ws, $cell, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.ClearContents()
}

function Append-Kayit {
    param($ws, [int]$newRow)

    $cellA = $ws.Range("A$newRow")
    $cellB = $ws.Range("B$newRow")
    $cellC = $ws.Range("C$newRow")
    $cellD = $ws.Range("D$newRow")
    $cellE = $ws.Range("E$newRow")
    $cellF = $ws.Range("F$newRow")

    Set-TextValue $ws $cellA $kayitNo
    Set-TextValue $ws $cellB $tarih
    Set-TextValue $ws $cellC $birim
    Set-TextValue $ws $cellD $parsel
    Set-TextValue $ws $cellE $is
    Set-TextValue $ws $cellF $personel
}

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Append-Kayit $wsKayitlar 25

$wsErdemli = $wb.Worksheets.Item("Erdemli")
Append-Kayit $wsErdemli 24
